# admin can remove staff, make middlename optional
#
# Update the "amount" column (B) sample data so each value is truncated to
# 3 digits, and flag every bordered row in the table (A1:E10) with a white
# fill so validation/error rows stand out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Trim the stray digit off the numbers in column B ------------------
$ws.Range("B1").Value = 123
$ws.Range("B2").Value = 432
$ws.Range("B3").Value = 876
$ws.Range("B4").Value = 543
$ws.Range("B6").Value = 147
$ws.Range("B8").Value = 309
$ws.Range("B9").Value = 294

# --- 2. Give the whole (bordered) table a white background fill -----------
$table = $ws.Range("A1:E10")
$table.Interior.ColorIndex = 2
